$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 32051.406
$ws.Range("I28").Value = 34105.168
$ws.Range("J28").Value = 1245
$ws.Range("K28").Value = 34105.168
$ws.Range("L28").Value = 1245
$ws.Range("M28").Value = -33620.168
$ws.Range("N28").Value = -2215
$ws.Range("H32").Value = 6166.3335
$ws.Range("I32").Value = 1833.3334
$ws.Range("K32").Value = 1833.3334
$ws.Range("M32").Value = -1507.3334
$ws.Range("H86").Value = 6935080.5
$ws.Range("I86").Value = 2154.111
$ws.Range("J86").Value = 10054897
$ws.Range("K86").Value = 2154.111
$ws.Range("L86").Value = 10054897
$ws.Range("M86").Value = -1031.111
$ws.Range("N86").Value = -10057143
$ws.Range("H89").Value = 6935080.5
$ws.Range("I89").Value = 2154.111
$ws.Range("J89").Value = 10054897
$ws.Range("K89").Value = 10770.555
$ws.Range("L89").Value = 50274485
$ws.Range("M89").Value = -5154.555
$ws.Range("N89").Value = -50285717
$ws.Range("H112").Value = 126754.875
$ws.Range("J112").Value = 2005.5714
$ws.Range("L112").Value = 6016.7142
$ws.Range("N112").Value = -8232.7142
$ws.Range("H138").Value = 2726.5088
$ws.Range("J138").Value = 2673
$ws.Range("L138").Value = 8019
$ws.Range("N138").Value = -18299

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 5003499.5
$ws.Range("I8").Value = 6669333.5
$ws.Range("K8").Value = 6669333.5
$ws.Range("M8").Value = -6669189.5
$ws.Range("H14").Value = 664
$ws.Range("I14").Value = 330
$ws.Range("K14").Value = 330
$ws.Range("M14").Value = -155
$ws.Range("H32").Value = 13087.419
$ws.Range("I32").Value = 8161.125
$ws.Range("J32").Value = 28413.666
$ws.Range("K32").Value = 8161.125
$ws.Range("L32").Value = 28413.666
$ws.Range("M32").Value = -7874.125
$ws.Range("N32").Value = -28987.666
$ws.Range("H61").Value = 3717.7856
$ws.Range("I61").Value = 2695.4546
$ws.Range("J61").Value = 7466.3335
$ws.Range("K61").Value = 2695.4546
$ws.Range("L61").Value = 7466.3335
$ws.Range("M61").Value = -2483.4546
$ws.Range("N61").Value = -7890.3335
$ws.Range("H74").Value = 38464790
$ws.Range("J74").Value = 4328.5
$ws.Range("L74").Value = 4328.5
$ws.Range("N74").Value = -6076.5
$ws.Range("H77").Value = 38464790
$ws.Range("J77").Value = 4328.5
$ws.Range("L77").Value = 21642.5
$ws.Range("N77").Value = -30378.5
$ws.Range("H122").Value = 4376
$ws.Range("I122").Value = 4107.875
$ws.Range("K122").Value = 12323.625
$ws.Range("M122").Value = -9873.625
$ws.Range("H132").Value = 2826.4546
$ws.Range("I132").Value = 2327.4333
$ws.Range("K132").Value = 6982.2999
$ws.Range("M132").Value = -4452.2999
$ws.Range("H136").Value = 3717.7856
$ws.Range("I136").Value = 2695.4546
$ws.Range("J136").Value = 7466.3335
$ws.Range("K136").Value = 8086.3638
$ws.Range("L136").Value = 22399.0005
$ws.Range("M136").Value = -5536.3638
$ws.Range("N136").Value = -27499.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 4000
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H16").Value = 6501.25
$ws.Range("I16").Value = 6002
$ws.Range("K16").Value = 6002
$ws.Range("M16").Value = -5832
$ws.Range("H105").Value = 2244.389
$ws.Range("I105").Value = 2093.3667
$ws.Range("K105").Value = 2093.3667
$ws.Range("M105").Value = -346.3667
$ws.Range("H134").Value = 2350.88
$ws.Range("I134").Value = 1777.289
$ws.Range("J134").Value = 7513.2
$ws.Range("K134").Value = 5331.867
$ws.Range("L134").Value = 22539.6
$ws.Range("M134").Value = -2796.867
$ws.Range("N134").Value = -27609.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 95.545456
$ws.Range("I7").Value = 87.625
$ws.Range("J7").Value = 116.666664
$ws.Range("K7").Value = 87.625
$ws.Range("L7").Value = 116.666664
$ws.Range("M7").Value = 25.375
$ws.Range("N7").Value = -342.666664
$ws.Range("H22").Value = 932.9375
$ws.Range("J22").Value = 1222.4445
$ws.Range("L22").Value = 1222.4445
$ws.Range("N22").Value = -1922.4445
$ws.Range("H70").Value = 23750
$ws.Range("J70").Value = 23750
$ws.Range("L70").Value = 23750
$ws.Range("N70").Value = -24380
$ws.Range("H73").Value = 23750
$ws.Range("J73").Value = 23750
$ws.Range("L73").Value = 23750
$ws.Range("N73").Value = -25934
$ws.Range("H134").Value = 4729.657
$ws.Range("I134").Value = 4205.5356
$ws.Range("K134").Value = 12616.6068
$ws.Range("M134").Value = -10081.6068
$ws.Range("H141").Value = 198986.88
$ws.Range("J141").Value = 198986.88
$ws.Range("L141").Value = 198986.88
$ws.Range("N141").Value = -209346.88

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 145
$ws.Range("J12").Value = 175
$ws.Range("L12").Value = 525
$ws.Range("N12").Value = -871
$ws.Range("H34").Value = 1888.4117
$ws.Range("J34").Value = 2199.1428
$ws.Range("L34").Value = 6597.428400000001
$ws.Range("N34").Value = -6765.428400000001
$ws.Range("H39").Value = 4620.75
$ws.Range("J39").Value = 3548.158
$ws.Range("L39").Value = 10644.474
$ws.Range("N39").Value = -11232.474
$ws.Range("H55").Value = 3642.75
$ws.Range("J55").Value = 3642.75
$ws.Range("L55").Value = 10928.25
$ws.Range("N55").Value = -11282.25
$ws.Range("H127").Value = 900
$ws.Range("J127").Value = 900
$ws.Range("L127").Value = 2700
$ws.Range("N127").Value = -12620
$ws.Range("H132").Value = 2792.491
$ws.Range("J132").Value = 2954.2173
$ws.Range("L132").Value = 26587.9557
$ws.Range("N132").Value = -31647.9557
$ws.Range("H137").Value = 2898.2593
$ws.Range("J137").Value = 3389
$ws.Range("L137").Value = 10167
$ws.Range("N137").Value = -20367

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3783.818
$ws.Range("I113").Value = 2833
$ws.Range("J113").Value = 4140.375
$ws.Range("K113").Value = 2833
$ws.Range("L113").Value = 4140.375
$ws.Range("M113").Value = -663
$ws.Range("N113").Value = -8480.375
$ws.Range("H132").Value = 3188.2163
$ws.Range("I132").Value = 2999.0557
$ws.Range("K132").Value = 8997.167099999999
$ws.Range("M132").Value = -6467.167099999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 566.4
$ws.Range("I9").Value = 705
$ws.Range("J9").Value = 474
$ws.Range("K9").Value = 705
$ws.Range("L9").Value = 474
$ws.Range("M9").Value = -481
$ws.Range("N9").Value = -922
$ws.Range("H132").Value = 18523502
$ws.Range("I132").Value = 25644218
$ws.Range("J132").Value = 9643.866
$ws.Range("K132").Value = 76932654
$ws.Range("L132").Value = 28931.598
$ws.Range("M132").Value = -76930124
$ws.Range("N132").Value = -33991.598
$ws.Range("H136").Value = 6471.787
$ws.Range("I136").Value = 6506.619
$ws.Range("J136").Value = 6179.2
$ws.Range("K136").Value = 19519.857
$ws.Range("L136").Value = 18537.6
$ws.Range("M136").Value = -16969.857
$ws.Range("N136").Value = -23637.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1450.75
$ws.Range("I6").Value = 4000
$ws.Range("J6").Value = 1086.5714
$ws.Range("K6").Value = 4000
$ws.Range("L6").Value = 1086.5714
$ws.Range("M6").Value = -3885
$ws.Range("N6").Value = -1316.5714
$ws.Range("H81").Value = 10006994
$ws.Range("J81").Value = 22232112
$ws.Range("L81").Value = 44464224
$ws.Range("N81").Value = -44466346
$ws.Range("H84").Value = 10006994
$ws.Range("J84").Value = 22232112
$ws.Range("L84").Value = 222321120
$ws.Range("N84").Value = -222331728
$ws.Range("H107").Value = 1627.9375
$ws.Range("I107").Value = 1289
$ws.Range("K107").Value = 3867
$ws.Range("M107").Value = -1947
$ws.Range("H132").Value = 18489.346
$ws.Range("I132").Value = 7480.737
$ws.Range("J132").Value = 48369.855
$ws.Range("K132").Value = 22442.211
$ws.Range("L132").Value = 145109.565
$ws.Range("M132").Value = -19912.211
$ws.Range("N132").Value = -150169.565
$ws.Range("H136").Value = 5027.476
$ws.Range("I136").Value = 6461.9546
$ws.Range("J136").Value = 3449.55
$ws.Range("K136").Value = 19385.8638
$ws.Range("L136").Value = 10348.65
$ws.Range("M136").Value = -16835.8638
$ws.Range("N136").Value = -15448.65
